# Commit: "Tue, Apr 07, 2020 12:07:09 PM"
#
# The authored change swaps the presentation's design theme (ppt/theme/theme1.xml,
# used by the slide master / the deck's visible design) from the "Integral" /
# "Red Violet" look to the plain default "Office Theme" / "Office" color
# palette (and, symmetrically, moves the old "Integral" colors onto the
# notes-master theme part, ppt/theme/theme2.xml). The font scheme and the
# fill/line/effect format scheme are byte-identical between the two theme
# parts already (both are the stock Arial / Office format scheme), so the
# only substantive, visible difference is the 12-slot theme color palette
# (and the cosmetic <a:theme>/<a:clrScheme> name= labels, which PowerPoint's
# object model does not expose a writer for - renaming a theme/color scheme
# isn't possible through Design/Theme/ThemeColorScheme from automation).
#
# Apply the new "Office" color scheme to the presentation's theme via the
# slide master's ThemeColorScheme, one RGB slot at a time, in the fixed
# COM order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      44546A  (COM RGB is 0xBBGGRR)
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink 954F72

# Best-effort: try to also relabel the color scheme the way the authored
# edit does (Red Violet -> Office). PowerPoint automation doesn't expose a
# way to rename a theme/color scheme, so this is a harmless no-op if the
# host ignores it. (Deliberately NOT touching Designs.Item(1).Name - on
# this host that writes to the slide master's <p:cSld name="..."> instead
# of the theme, which is not part of the intended change.)
try { $cs.Name = "Office" } catch {}
